# CS133JS Lab07 Rubric - "Updated rubrics for total score of 40"
#
# The possible-points column (and the mirrored score column on the Score
# sheet) is rescored so the rubric totals 40 instead of 50. Row 7 ("Part 1"
# header) never carried points and its placeholder D/E cells are cleared
# entirely. The now-unused "Possible"/"Score" cells outside the graded rows
# are also cleared down to nothing (no value, no format) to match a clean
# re-entry of the column.

$wb = $excel.ActiveWorkbook
$wsRubric = $wb.Worksheets.Item("Rubric")
$wsScore  = $wb.Worksheets.Item("Score")

# New "Possible points" per criteria row (row -> value), row 7 has none.
$newPoints = @{
    6  = 8
    8  = 2
    9  = 4
    10 = 3
    11 = 1
    12 = 3
    13 = 5
    14 = 3
    15 = 1
    16 = 3
    17 = 2
    18 = 3
    19 = 2
}

foreach ($sheet in @($wsRubric, $wsScore)) {
    # Wipe the old values/formatting for the whole editable block first
    # (Clear removes both content and the applied "D/E points" style, so
    # rows with no points left behind no placeholder cell at all).
    $sheet.Range("D6:E19").Clear()

    foreach ($row in $newPoints.Keys) {
        $value = $newPoints[$row]
        $sheet.Cells.Item($row, 4).Value = $value   # column D - Possible
    }
}

# The Score sheet also mirrors the same numbers into column E (the actual
# score achieved, which in this rubric template starts equal to Possible).
foreach ($row in $newPoints.Keys) {
    $value = $newPoints[$row]
    $wsScore.Cells.Item($row, 5).Value = $value     # column E - Score
}

# Recalculate the Total rows (D21 / E21 already hold SUM formulas).
$wb.Application.Calculate()

# Restore the view/selection state recorded in the saved workbook: Score
# sheet shows D6:D19 selected, Rubric sheet is the active tab with F12
# selected.
$wsScore.Select()
$wsScore.Range("D6:D19").Select()

$wsRubric.Select()
$wsRubric.Range("F12").Select()
